$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.741.05"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.300.43"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.71"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "574.55"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.874.40"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.42"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.099.86"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.293.06"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "441.74"
$ws.Range("E18").Value = "  +9.92%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.76"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.26"
$ws.Range("E22").Value = "  +4.53%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.513"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.432.88"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.81"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.33"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.77"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.77"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.54"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.783"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.726.83"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.26"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.17"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0671"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.73"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "327.04"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("E51").Value = "  -0.74%  "
